$d = $word.ActiveDocument
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

# 1) "Tipo" cell in the first table is split across two runs ("P" + "rimário, essencial").
#    Re-typing the same visible text over the whole match merges it into a single run.
#    Only the first occurrence (Replace=1/wdReplaceOne) is touched; the second table's
#    "Tipo" cell already has the text in a single run and must stay untouched.
$find.Execute("Primário, essencial", $false, $false, $false, $false, $false, $true, 1, $false, "Primário, essencial", 1) | Out-Null

# 2) The sentence "...tendo que ser retirado antes que ele feche." is split across three runs
#    (main text / "ele" / " feche."). Re-typing the same text merges them into a single run.
$find.Execute("tendo que ser retirado antes que ele feche.", $false, $false, $false, $false, $false, $true, 1, $false, "tendo que ser retirado antes que ele feche.", 1) | Out-Null

# 3) Fix grammar typo: "não fora retirado" -> "não for retirado"
$find.Execute("não fora retirado", $false, $false, $false, $false, $false, $true, 1, $false, "não for retirado", 1) | Out-Null
